$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlink annotations up front (cell text values are untouched);
# they will be re-created below, after the row shift, so that each hyperlink's target
# stays correctly paired with the URL text actually displayed in its cell.
$ws.Range("F2").Hyperlinks.Delete()

# Insert a new row at position 13; this shifts existing rows 13-19 down to 14-20
# (carrying along their values, number formats and styles).
$ws.Rows.Item(13).Insert()

$newTimestamp = "2025-11-11 01:50:31"

# Refresh the "取得日時" timestamp (column A) on every data row (2-20)
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Populate the brand-new row 13 with the newly scraped listing
$ws.Cells.Item(13, 2).Value = "【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5431322"
$ws.Cells.Item(13, 7).Value = 25

# Re-create the URL hyperlinks for every data row, top to bottom, so relationship
# ids line up with final row order and every link target matches its cell text.
$urls = @{
    2 = "https://www.lancers.jp/work/detail/5430365"
    3 = "https://www.lancers.jp/work/detail/5431299"
    4 = "https://www.lancers.jp/work/detail/5430993"
    5 = "https://www.lancers.jp/work/detail/5430799"
    6 = "https://www.lancers.jp/work/detail/5430954"
    7 = "https://www.lancers.jp/work/detail/5416665"
    8 = "https://www.lancers.jp/work/detail/5430436"
    9 = "https://www.lancers.jp/work/detail/5431107"
    10 = "https://www.lancers.jp/work/detail/5431051"
    11 = "https://www.lancers.jp/work/detail/5430951"
    12 = "https://www.lancers.jp/work/detail/5418064"
    13 = "https://www.lancers.jp/work/detail/5431322"
    14 = "https://www.lancers.jp/work/detail/5431085"
    15 = "https://www.lancers.jp/work/detail/5428756"
    16 = "https://www.lancers.jp/work/detail/5428755"
    17 = "https://www.lancers.jp/work/detail/5431069"
    18 = "https://www.lancers.jp/work/detail/5431036"
    19 = "https://www.lancers.jp/work/detail/5431276"
    20 = "https://www.lancers.jp/work/detail/5417622"
}
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
}
$ws.Range("F2:F20").Style = "Hyperlink"
